# correction in sa algorithm and 746 logs
# Update the "Fitness" (column C) values for rows 2-98 on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5 -> 7672
$ws.Range("C2:C5").Value = 7672

# Rows 6-25 -> 7312
$ws.Range("C6:C25").Value = 7312

# Rows 26-98 -> 7310
$ws.Range("C26:C98").Value = 7310
